# Edit: rename "Gen" header to "MaxFES", convert column A generation counts
# into normalized MaxFES fractions, drop the "Run 50" column (old column AZ)
# and recompute the "Mean" column (now shifted into AZ) over Run 0..Run 49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the "Run 50" column (column AZ, index 52). This shifts the old
#    "Mean" column (was BA) one position left, into AZ.
$ws.Columns.Item(52).Delete()

# 2. Header rename: "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# 3. New column-A values (normalized MaxFES fractions) for rows 2-14
$newA = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $newA.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value2 = $newA[$i]
}

# 4. Recompute "Mean" (now in column AZ / 52) over Run 0..Run 49 (cols B..AY = 2..51)
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 52)
    $cell.Formula = "=AVERAGE(B" + $r + ":AY" + $r + ")"
    $cell.Value2 = $cell.Value2
}
